# Update the "Written" page count for the day of 2021-10-09 (row 20, col I)
# from 0 to 5.5 pages. All the other cells in the sheet (B5, B11, B12, D11,
# I31, F3, F4) are formulas that depend on this value (directly or
# transitively) and will recalculate automatically. F2 depends on TODAY()
# and will also recalculate against the workbook's pinned clock.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(20, 9).Value = 5.5

# Force a full recalculation so every dependent formula's cached value is
# refreshed (matches Excel's automatic-calculation behavior).
$excel.Calculate()

# Move the active selection to I21, matching the author's final cursor
# position after entering the new value.
$ws.Range("I21").Select()
